# Decile scaling to improve contrast
# Adds Min / Decile / Max summary block (columns H:L, rows 1-13) to the
# weather_data sheet, mirroring the precip_day / precip_avg / temp_avg
# columns (D:F) with MIN / PERCENTILE / MAX formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - mirrored headers (precip_day / precip_avg / temp_avg)
# ---------------------------------------------------------------------
$ws.Range("J1").Formula = '=D1'
$ws.Range("K1").Formula = '=E1'
$ws.Range("L1").Formula = '=F1'
$ws.Range("J1:L1").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Row 2 - Min row
# ---------------------------------------------------------------------
$ws.Range("I2").Value = "Min"

$ws.Range("J2").Formula = '=MIN(D2:D161)'
$ws.Range("K2").Formula = '=MIN(E2:E161)'
$ws.Range("J2:K2").NumberFormat = "0.00"

$ws.Range("L2").Formula = '=MIN(F2:F161)'

# ---------------------------------------------------------------------
# Rows 3-12 - Decile rows (10% .. 100%)
# ---------------------------------------------------------------------
$ws.Range("H3").Value = "Decile"
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H3").VerticalAlignment = -4108

$deciles = @(10, 20, 30, 40, 50, 60, 70, 80, 90, 100)
for ($i = 0; $i -lt $deciles.Length; $i++) {
    $row = 3 + $i
    $pct = $deciles[$i] / 100

    $ws.Cells.Item($row, 9).Value = $deciles[$i]
    $ws.Cells.Item($row, 9).NumberFormat = "0.00"

    if ($row -eq 3) {
        $ws.Range("J3").Formula = '=PERCENTILE(D$2:D$161,' + $pct + ')'
    } else {
        $ws.Cells.Item($row, 10).Formula = '=PERCENTILE($D$2:$D$161,' + $pct + ')'
    }
    $ws.Cells.Item($row, 10).NumberFormat = "0.0000"

    $ws.Cells.Item($row, 11).Formula = '=PERCENTILE(E$2:E$161,' + $pct + ')'

    $ws.Cells.Item($row, 12).Formula = '=PERCENTILE(F$2:F$161,' + $pct + ')'
}

# Merge the "Decile" label down the full decile block (style already set
# on H3 above so the merge just propagates it, avoiding extra style churn)
$ws.Range("H3:H12").Merge()

# ---------------------------------------------------------------------
# Row 13 - Max row
# ---------------------------------------------------------------------
$ws.Range("I13").Value = "Max"

$ws.Range("J13").Formula = '=MAX(D2:D161)'
$ws.Range("J13").NumberFormat = "0.00"

$ws.Range("K13").Formula = '=MAX(E2:E161)'
$ws.Range("K13").NumberFormat = "0.00000"

$ws.Range("L13").Formula = '=MAX(F2:F161)'

# Apply the high-precision temp_avg format last, across the whole L block
# (Min / deciles / Max), matching the single shared numFmt used there.
$ws.Range("L2").NumberFormat = "0.00000000000"
$ws.Range("L3:L12").NumberFormat = "0.00000000000"
$ws.Range("L13").NumberFormat = "0.00000000000"

# ---------------------------------------------------------------------
# Cosmetics: column widths + selection, matching the authored view
# ---------------------------------------------------------------------
$ws.Columns("J").ColumnWidth = 9.5
$ws.Columns("K").ColumnWidth = 12.1640625
$ws.Columns("L").ColumnWidth = 14.6640625

$ws.Range("L2:L13").Select()
